# FlaglerEdits_v2: "updated 'contact us' link in N01P"
#
# 1) Collapse the split hyperlink-text runs (artifacts of earlier
#    per-character edits) back into a single run per hyperlink, for
#    every URL that had been split across multiple <w:r> elements.
# 2) Highlight (green) the "Link edit - Near the base..." question
#    in the N1P section and add a new paragraph right after it
#    answering where the "contact us" link now points.
# 3) Drop the stale lastRenderedPageBreak cached on the "No edits."
#    paragraph following the LP hyperlink.

$d = $word.ActiveDocument

# --- 1) Re-join hyperlink text that had been split across runs ---------
$hyperlinkUrls = @(
    "https://liammoor.seeflagler.com/E1.html",
    "https://liammoor.seeflagler.com/E2.html",
    "https://liammoor.seeflagler.com/E3.html",
    "https://liammoor.seeflagler.com/E4.html",
    "https://liammoor.seeflagler.com/E5.html",
    "https://liammoor.seeflagler.com/E6.html",
    "https://liammoor.seeflagler.com/E7.html",
    "https://liammoor.seeflagler.com/E8.html",
    "https://liammoor.seeflagler.com/E9.html",
    "https://liammoor.seeflagler.com/E10.html",
    "https://liammoor.seeflagler.com/E11.html",
    "https://liammoor.seeflagler.com/E12.html",
    "https://seeflagler2023.secure.mdl.io/E2.html",
    "https://seeflagler2023.secure.mdl.io/E3.html",
    "https://seeflagler2023.secure.mdl.io/E4.html",
    "https://seeflagler2023.secure.mdl.io/E5.html",
    "https://seeflagler2023.secure.mdl.io/E9.html",
    "https://seeflagler2023.secure.mdl.io/N1.html"
)

foreach ($url in $hyperlinkUrls) {
    $d.Content.Find.Execute($url, $true, $false, $false, $false, $false,
                             $true, 1, $false, $url, 2) | Out-Null
}

# --- 2) "contact us" link follow-up in N1P ------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Text = "Near the base of the email there is a link option to"
$rng.Find.Forward = $true
$rng.Find.Wrap = 1
$rng.Find.Execute() | Out-Null

$questionPara = $rng.Paragraphs(1)
$questionPara.Range.HighlightColorIndex = 4  # wdBrightGreen -> w:highlight val="green"

$questionPara.Range.InsertParagraphAfter()
$answerPara = $questionPara.Next()
$answerPara.Range.Text = "Now linked to flagler.edu/admissions--aid/contact us/"
$answerPara.Range.HighlightColorIndex = 0    # wdNoHighlight, keep the new line unhighlighted

# --- 3) Remove the stale lastRenderedPageBreak on "No edits." (after LP) -
$lpRng = $d.Content
$lpRng.Find.ClearFormatting()
$lpRng.Find.Text = "lp.html"
$lpRng.Find.Forward = $true
$lpRng.Find.Wrap = 1
$lpRng.Find.Execute() | Out-Null

$lpPara = $lpRng.Paragraphs(1)
$noEditsPara = $lpPara.Next()
$noEditsPara.Range.Text = "No edits."
